$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.458.54"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.878.11"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  -2.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5112"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3944"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08408"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.109"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.73"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.262"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "1.875.49"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.273"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.014"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001106"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.17"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06742"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.69"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.952"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("D23").Value = "28.514.35"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.267"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.97%  "
$ws.Range("D26").Value = "2.093.14"
$ws.Range("E26").Value = "  -1.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.49"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.72"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.382"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.41"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.792"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.629"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02429"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06500"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2188"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.924"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.267"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.191"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.085"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6442"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.18"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6061"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.09"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.711"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.014"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.201"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.30%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.78%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.208"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.08%  "
